$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 350-351 (everything from the old row 350 onward
# shifts down by two rows, e.g. old row 350 -> new row 352, old row 481 ->
# new row 483). The workbook dimension grows from A1:R481 to A1:R483.
$ws.Rows("350:351").Insert()

# Populate the newly inserted row 350 ("Primera" quality) with the new
# weekly price data.
$ws.Range("A350").Value = 8
$ws.Range("B350").Value = "Terminal La Palmera de La Serena"
$ws.Range("C350").Value = "Coquimbo"
$ws.Range("D350").Value = 45120
$ws.Range("E350").Value = 4
$ws.Range("F350").Value = 100114014
$ws.Range("G350").Value = "Betarraga"
$ws.Range("H350").Value = "Sin especificar"
$ws.Range("I350").Value = "Primera"
$ws.Range("J350").Value = 1960
$ws.Range("K350").Value = 550
$ws.Range("L350").Value = 600
$ws.Range("M350").Value = 575
$ws.Range("N350").Value = "`$/paquete 3 unidades"
$ws.Range("O350").Value = "Provincia del Elquí"
$ws.Range("P350").Value = 192
$ws.Range("Q350").Value = 3
$ws.Range("R350").Value = "Hortaliza"

# Populate the newly inserted row 351 ("Segunda" quality) with the new
# weekly price data.
$ws.Range("A351").Value = 8
$ws.Range("B351").Value = "Terminal La Palmera de La Serena"
$ws.Range("C351").Value = "Coquimbo"
$ws.Range("D351").Value = 45120
$ws.Range("E351").Value = 4
$ws.Range("F351").Value = 100114014
$ws.Range("G351").Value = "Betarraga"
$ws.Range("H351").Value = "Sin especificar"
$ws.Range("I351").Value = "Segunda"
$ws.Range("J351").Value = 1560
$ws.Range("K351").Value = 450
$ws.Range("L351").Value = 500
$ws.Range("M351").Value = 475
$ws.Range("N351").Value = "`$/paquete 3 unidades"
$ws.Range("O351").Value = "Provincia del Elquí"
$ws.Range("P351").Value = 158
$ws.Range("Q351").Value = 3
$ws.Range("R351").Value = "Hortaliza"
